# Append a new "5. Docker" section at the end of the document, after the
# paragraph that ends with "...make the bool true."
#
# The new section consists of:
#   Heading2 : "5. Docker"
#   Heading3 : "docker-compose.yml"
#   Normal   : "- Copy old game's line in docker-compose.yml"
#   Normal   : "- Change main namespace"
#   Normal   : "- Change "context" to reference new game file"
#   Normal   : "- Change "container_name" to "game-(gamename)""
#   Normal   : "- Change port Number (1950 for BucKart)"

$d = $word.ActiveDocument

# Locate the very last paragraph in the document (the one ending in "true.")
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)

# Create a brand-new, empty paragraph right after it to host the inserted
# content (keeps the original paragraph's text/formatting untouched).
$insertionRange = $lastPara.Range
$insertionRange.Collapse(0)
$null = $insertionRange.InsertParagraphAfter()

$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newIndex)

# Build the raw WordprocessingML for the seven new paragraphs and drop it
# in as the new paragraph's content, preserving heading styles, spell-check
# markers and the page-break hint exactly as authored.
$xml = @'
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>5</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Docker</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>docker-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>compose.yml</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">- </w:t></w:r><w:r><w:t xml:space="preserve">Copy old game’s line in </w:t></w:r><w:r><w:t>docker-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>compose.yml</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>- Change main namespace</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">- Change “context” </w:t></w:r><w:r><w:t>to reference new game file</w:t></w:r></w:p><w:p><w:r><w:t>- Change “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>container_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve"> to “</w:t></w:r><w:r><w:t>game-(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gamename</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)”</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">- Change port Number (1950 for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BucKart</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@

$null = $newPara.Range.InsertXML($xml)
